$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "2023-12-06 09:46:54"
$ws.Range("B12").Value = 0.0016

$ws.Range("A13").Value = "2023-12-06 09:48:40"
$ws.Range("B13").Value = 0.006200000000000002
